$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 70 (constant columns copied from the existing pattern)
$ws.Range("A70").Value = 9
$ws.Range("B70").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C70").Value = "Metropolitana"
$ws.Range("D70").Value = 45013
$ws.Range("D70").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E70").Value = 13
$ws.Range("F70").Value = 100112010
$ws.Range("G70").Value = "Achicoria"
$ws.Range("H70").Value = "Sin especificar"
$ws.Range("I70").Value = "Primera"
$ws.Range("J70").Value = 70
$ws.Range("K70").Value = 7000
$ws.Range("L70").Value = 7000
$ws.Range("M70").Value = 7000
$ws.Range("N70").Value = "$/caja 16 unidades"
$ws.Range("O70").Value = "Provincia de Quillota"
$ws.Range("P70").Value = 438
$ws.Range("Q70").Value = 16
$ws.Range("R70").Value = "Hortaliza"

# Shift weekly records: rows 24-69 take the prior rows D/J/K/L/M/O/P values;
# row 23 becomes the newest record (the row 70 data above was the old row 69).
$ws.Range("D69").Value = 44999
$ws.Range("J69").Value = 70
$ws.Range("K69").Value = 7000
$ws.Range("L69").Value = 7000
$ws.Range("M69").Value = 7000
$ws.Range("O69").Value = "Provincia de Quillota"
$ws.Range("P69").Value = 438
$ws.Range("D68").Value = 44467
$ws.Range("J68").Value = 52
$ws.Range("K68").Value = 5000
$ws.Range("L68").Value = 6000
$ws.Range("M68").Value = 5500
$ws.Range("O68").Value = "Provincia de Quillota"
$ws.Range("P68").Value = 344
$ws.Range("D67").Value = 45069
$ws.Range("J67").Value = 70
$ws.Range("K67").Value = 7000
$ws.Range("L67").Value = 7000
$ws.Range("M67").Value = 7000
$ws.Range("O67").Value = "Provincia de Quillota"
$ws.Range("P67").Value = 438
$ws.Range("D66").Value = 45085
$ws.Range("J66").Value = 90
$ws.Range("K66").Value = 7000
$ws.Range("L66").Value = 7000
$ws.Range("M66").Value = 7000
$ws.Range("O66").Value = "Provincia de Quillota"
$ws.Range("P66").Value = 438
$ws.Range("D65").Value = 44967
$ws.Range("J65").Value = 70
$ws.Range("K65").Value = 7000
$ws.Range("L65").Value = 7000
$ws.Range("M65").Value = 7000
$ws.Range("O65").Value = "Provincia de Quillota"
$ws.Range("P65").Value = 438
$ws.Range("D64").Value = 44474
$ws.Range("J64").Value = 52
$ws.Range("K64").Value = 5000
$ws.Range("L64").Value = 6000
$ws.Range("M64").Value = 5500
$ws.Range("O64").Value = "Provincia de Quillota"
$ws.Range("P64").Value = 344
$ws.Range("D63").Value = 44782
$ws.Range("J63").Value = 70
$ws.Range("K63").Value = 6000
$ws.Range("L63").Value = 6000
$ws.Range("M63").Value = 6000
$ws.Range("O63").Value = "Región Metropolitana"
$ws.Range("P63").Value = 375
$ws.Range("D62").Value = 44691
$ws.Range("J62").Value = 61
$ws.Range("K62").Value = 6000
$ws.Range("L62").Value = 7000
$ws.Range("M62").Value = 6508
$ws.Range("O62").Value = "Provincia de Quillota"
$ws.Range("P62").Value = 407
$ws.Range("D61").Value = 44938
$ws.Range("J61").Value = 70
$ws.Range("K61").Value = 7000
$ws.Range("L61").Value = 7000
$ws.Range("M61").Value = 7000
$ws.Range("O61").Value = "Provincia de Quillota"
$ws.Range("P61").Value = 438
$ws.Range("D60").Value = 45055
$ws.Range("J60").Value = 70
$ws.Range("K60").Value = 7000
$ws.Range("L60").Value = 7000
$ws.Range("M60").Value = 7000
$ws.Range("O60").Value = "Provincia de Quillota"
$ws.Range("P60").Value = 438
$ws.Range("D59").Value = 44946
$ws.Range("J59").Value = 70
$ws.Range("K59").Value = 7000
$ws.Range("L59").Value = 7000
$ws.Range("M59").Value = 7000
$ws.Range("O59").Value = "Provincia de Quillota"
$ws.Range("P59").Value = 438
$ws.Range("D58").Value = 44698
$ws.Range("J58").Value = 34
$ws.Range("K58").Value = 6000
$ws.Range("L58").Value = 7000
$ws.Range("M58").Value = 6500
$ws.Range("O58").Value = "Provincia de Quillota"
$ws.Range("P58").Value = 406
$ws.Range("D57").Value = 44957
$ws.Range("J57").Value = 70
$ws.Range("K57").Value = 7000
$ws.Range("L57").Value = 7000
$ws.Range("M57").Value = 7000
$ws.Range("O57").Value = "Provincia de Quillota"
$ws.Range("P57").Value = 438
$ws.Range("D56").Value = 44910
$ws.Range("J56").Value = 70
$ws.Range("K56").Value = 6000
$ws.Range("L56").Value = 7000
$ws.Range("M56").Value = 6500
$ws.Range("O56").Value = "Provincia de Quillota"
$ws.Range("P56").Value = 406
$ws.Range("D55").Value = 44350
$ws.Range("J55").Value = 25
$ws.Range("K55").Value = 6000
$ws.Range("L55").Value = 6000
$ws.Range("M55").Value = 6000
$ws.Range("O55").Value = "Provincia de Quillota"
$ws.Range("P55").Value = 375
$ws.Range("D54").Value = 45062
$ws.Range("J54").Value = 52
$ws.Range("K54").Value = 7000
$ws.Range("L54").Value = 7000
$ws.Range("M54").Value = 7000
$ws.Range("O54").Value = "Provincia de Quillota"
$ws.Range("P54").Value = 438
$ws.Range("D53").Value = 44341
$ws.Range("J53").Value = 51
$ws.Range("K53").Value = 5500
$ws.Range("L53").Value = 6000
$ws.Range("M53").Value = 5755
$ws.Range("O53").Value = "Provincia de Quillota"
$ws.Range("P53").Value = 360
$ws.Range("D52").Value = 44953
$ws.Range("J52").Value = 90
$ws.Range("K52").Value = 7000
$ws.Range("L52").Value = 7000
$ws.Range("M52").Value = 7000
$ws.Range("O52").Value = "Provincia de Quillota"
$ws.Range("P52").Value = 438
$ws.Range("D51").Value = 44306
$ws.Range("J51").Value = 50
$ws.Range("K51").Value = 6000
$ws.Range("L51").Value = 6000
$ws.Range("M51").Value = 6000
$ws.Range("O51").Value = "Provincia de Quillota"
$ws.Range("P51").Value = 375
$ws.Range("D50").Value = 45091
$ws.Range("J50").Value = 70
$ws.Range("K50").Value = 6000
$ws.Range("L50").Value = 6000
$ws.Range("M50").Value = 6000
$ws.Range("O50").Value = "Provincia de Quillota"
$ws.Range("P50").Value = 375
$ws.Range("D49").Value = 45097
$ws.Range("J49").Value = 70
$ws.Range("K49").Value = 7000
$ws.Range("L49").Value = 7000
$ws.Range("M49").Value = 7000
$ws.Range("O49").Value = "Provincia de Quillota"
$ws.Range("P49").Value = 438
$ws.Range("D48").Value = 45022
$ws.Range("J48").Value = 60
$ws.Range("K48").Value = 7000
$ws.Range("L48").Value = 7000
$ws.Range("M48").Value = 7000
$ws.Range("O48").Value = "Provincia de Quillota"
$ws.Range("P48").Value = 438
$ws.Range("D47").Value = 45083
$ws.Range("J47").Value = 70
$ws.Range("K47").Value = 7000
$ws.Range("L47").Value = 7000
$ws.Range("M47").Value = 7000
$ws.Range("O47").Value = "Provincia de Quillota"
$ws.Range("P47").Value = 438
$ws.Range("D46").Value = 44943
$ws.Range("J46").Value = 70
$ws.Range("K46").Value = 7000
$ws.Range("L46").Value = 7000
$ws.Range("M46").Value = 7000
$ws.Range("O46").Value = "Provincia de Quillota"
$ws.Range("P46").Value = 438
$ws.Range("D45").Value = 44932
$ws.Range("J45").Value = 70
$ws.Range("K45").Value = 6000
$ws.Range("L45").Value = 7000
$ws.Range("M45").Value = 6500
$ws.Range("O45").Value = "Provincia de Quillota"
$ws.Range("P45").Value = 406
$ws.Range("D44").Value = 45037
$ws.Range("J44").Value = 90
$ws.Range("K44").Value = 7000
$ws.Range("L44").Value = 7000
$ws.Range("M44").Value = 7000
$ws.Range("O44").Value = "Provincia de Quillota"
$ws.Range("P44").Value = 438
$ws.Range("D43").Value = 45015
$ws.Range("J43").Value = 90
$ws.Range("K43").Value = 7000
$ws.Range("L43").Value = 7000
$ws.Range("M43").Value = 7000
$ws.Range("O43").Value = "Provincia de Quillota"
$ws.Range("P43").Value = 438
$ws.Range("D42").Value = 45020
$ws.Range("J42").Value = 75
$ws.Range("K42").Value = 7000
$ws.Range("L42").Value = 7000
$ws.Range("M42").Value = 7000
$ws.Range("O42").Value = "Provincia de Quillota"
$ws.Range("P42").Value = 438
$ws.Range("D41").Value = 44582
$ws.Range("J41").Value = 52
$ws.Range("K41").Value = 7000
$ws.Range("L41").Value = 7000
$ws.Range("M41").Value = 7000
$ws.Range("O41").Value = "Provincia de Quillota"
$ws.Range("P41").Value = 438
$ws.Range("D40").Value = 44715
$ws.Range("J40").Value = 70
$ws.Range("K40").Value = 5000
$ws.Range("L40").Value = 6000
$ws.Range("M40").Value = 5500
$ws.Range("O40").Value = "Provincia de Quillota"
$ws.Range("P40").Value = 344
$ws.Range("D39").Value = 44573
$ws.Range("J39").Value = 34
$ws.Range("K39").Value = 8000
$ws.Range("L39").Value = 8000
$ws.Range("M39").Value = 8000
$ws.Range("O39").Value = "Provincia de Quillota"
$ws.Range("P39").Value = 500
$ws.Range("D38").Value = 44403
$ws.Range("J38").Value = 43
$ws.Range("K38").Value = 6000
$ws.Range("L38").Value = 6000
$ws.Range("M38").Value = 6000
$ws.Range("O38").Value = "Provincia de Quillota"
$ws.Range("P38").Value = 375
$ws.Range("D37").Value = 45030
$ws.Range("J37").Value = 90
$ws.Range("K37").Value = 7000
$ws.Range("L37").Value = 7000
$ws.Range("M37").Value = 7000
$ws.Range("O37").Value = "Provincia de Quillota"
$ws.Range("P37").Value = 438
$ws.Range("D36").Value = 44971
$ws.Range("J36").Value = 160
$ws.Range("K36").Value = 7000
$ws.Range("L36").Value = 7000
$ws.Range("M36").Value = 7000
$ws.Range("O36").Value = "Provincia de Quillota"
$ws.Range("P36").Value = 438
$ws.Range("D35").Value = 45007
$ws.Range("J35").Value = 70
$ws.Range("K35").Value = 7000
$ws.Range("L35").Value = 7000
$ws.Range("M35").Value = 7000
$ws.Range("O35").Value = "Provincia de Quillota"
$ws.Range("P35").Value = 438
$ws.Range("D34").Value = 45042
$ws.Range("J34").Value = 70
$ws.Range("K34").Value = 7000
$ws.Range("L34").Value = 7000
$ws.Range("M34").Value = 7000
$ws.Range("O34").Value = "Provincia de Quillota"
$ws.Range("P34").Value = 438
$ws.Range("D33").Value = 44455
$ws.Range("J33").Value = 52
$ws.Range("K33").Value = 5000
$ws.Range("L33").Value = 6000
$ws.Range("M33").Value = 5500
$ws.Range("O33").Value = "Provincia de Quillota"
$ws.Range("P33").Value = 344
$ws.Range("D32").Value = 45079
$ws.Range("J32").Value = 160
$ws.Range("K32").Value = 7000
$ws.Range("L32").Value = 7000
$ws.Range("M32").Value = 7000
$ws.Range("O32").Value = "Provincia de Quillota"
$ws.Range("P32").Value = 438
$ws.Range("D31").Value = 44407
$ws.Range("J31").Value = 45
$ws.Range("K31").Value = 5500
$ws.Range("L31").Value = 6000
$ws.Range("M31").Value = 5744
$ws.Range("O31").Value = "Provincia de Quillota"
$ws.Range("P31").Value = 359
$ws.Range("D30").Value = 44371
$ws.Range("J30").Value = 34
$ws.Range("K30").Value = 5500
$ws.Range("L30").Value = 6000
$ws.Range("M30").Value = 5750
$ws.Range("O30").Value = "Provincia de Quillota"
$ws.Range("P30").Value = 359
$ws.Range("D29").Value = 44575
$ws.Range("J29").Value = 61
$ws.Range("K29").Value = 8000
$ws.Range("L29").Value = 8000
$ws.Range("M29").Value = 8000
$ws.Range("O29").Value = "Provincia de Quillota"
$ws.Range("P29").Value = 500
$ws.Range("D28").Value = 44355
$ws.Range("J28").Value = 25
$ws.Range("K28").Value = 6000
$ws.Range("L28").Value = 6000
$ws.Range("M28").Value = 6000
$ws.Range("O28").Value = "Provincia de Quillota"
$ws.Range("P28").Value = 375
$ws.Range("D27").Value = 44994
$ws.Range("J27").Value = 70
$ws.Range("K27").Value = 7000
$ws.Range("L27").Value = 7000
$ws.Range("M27").Value = 7000
$ws.Range("O27").Value = "Provincia de Quillota"
$ws.Range("P27").Value = 438
$ws.Range("D26").Value = 44442
$ws.Range("J26").Value = 25
$ws.Range("K26").Value = 6000
$ws.Range("L26").Value = 7000
$ws.Range("M26").Value = 6480
$ws.Range("O26").Value = "Provincia de Quillota"
$ws.Range("P26").Value = 405
$ws.Range("D25").Value = 44589
$ws.Range("J25").Value = 52
$ws.Range("K25").Value = 8000
$ws.Range("L25").Value = 8000
$ws.Range("M25").Value = 8000
$ws.Range("O25").Value = "Provincia de Quillota"
$ws.Range("P25").Value = 500
$ws.Range("D24").Value = 44358
$ws.Range("J24").Value = 52
$ws.Range("K24").Value = 6000
$ws.Range("L24").Value = 6000
$ws.Range("M24").Value = 6000
$ws.Range("O24").Value = "Provincia de Quillota"
$ws.Range("P24").Value = 375
$ws.Range("D23").Value = 45100
$ws.Range("J23").Value = 52
$ws.Range("K23").Value = 8000
$ws.Range("L23").Value = 8000
$ws.Range("M23").Value = 8000
$ws.Range("O23").Value = "Provincia de Quillota"
$ws.Range("P23").Value = 500
